$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "ebike"
$ws.Range("B14").Value = 20

$ws.Range("B15").Select()
